# Updates the cryptos worksheet with the latest price/volume snapshot.
# Each cell is set individually to keep the text formatting (e.g. "60.041.82",
# "0.0000172") exactly as-is. Price cells in column D can look like numbers or
# dates to Excel's auto-detection, so we prefix them with a leading apostrophe
# (the same trick a human would use when typing into a cell) which forces
# Excel to store them as literal text without altering the displayed digits.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'59.973.15"
$ws.Range("E2").Value = "  +2.12%  "
$ws.Range("D3").Value = "'3.190.36"
$ws.Range("E3").Value = "  +1.13%  "
$ws.Range("D5").Value = "'536.17"
$ws.Range("E5").Value = "  +0.84%  "
$ws.Range("D6").Value = "'145.27"
$ws.Range("E6").Value = "  +3.62%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "'0.528"
$ws.Range("E8").Value = "  -1.73%  "
$ws.Range("D9").Value = "'7.32"
$ws.Range("E9").Value = "  -0.19%  "
$ws.Range("D10").Value = "'0.112"
$ws.Range("E10").Value = "  +2.15%  "
$ws.Range("E11").Value = "  -1.05%  "
$ws.Range("D12").Value = "'3.744.20"
$ws.Range("E12").Value = "  +1.17%  "
$ws.Range("E13").Value = "  -1.76%  "
$ws.Range("D14").Value = "'25.87"
$ws.Range("E14").Value = "  -1.16%  "
$ws.Range("D15").Value = "'0.0000172"
$ws.Range("E15").Value = "  +0.15%  "
$ws.Range("D16").Value = "'60.031.65"
$ws.Range("E16").Value = "  +2.14%  "
$ws.Range("D17").Value = "'3.193.78"
$ws.Range("E17").Value = "  +0.89%  "
$ws.Range("D18").Value = "'6.27"
$ws.Range("E18").Value = "  +0.14%  "
$ws.Range("D19").Value = "'13.27"
$ws.Range("E20").Value = "  +0.36%  "
$ws.Range("D21").Value = "'369.46"
$ws.Range("E21").Value = "  -0.99%  "
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("D24").Value = "'69.43"
$ws.Range("E24").Value = "  -1.04%  "
$ws.Range("E25").Value = "  +0.72%  "
$ws.Range("D26").Value = "'8.67"
$ws.Range("E26").Value = "  +5.29%  "
$ws.Range("D27").Value = "'1.01"
$ws.Range("E27").Value = "  +1.48%  "
$ws.Range("E28").Value = "  +0.27%  "
$ws.Range("D29").Value = "'22.42"
$ws.Range("E29").Value = "  +0.96%  "
$ws.Range("E30").Value = "  +0.43%  "
$ws.Range("E31").Value = "  -1.12%  "
$ws.Range("D32").Value = "'5.27"
$ws.Range("E32").Value = "  +1.38%  "
$ws.Range("E33").Value = "  +1.61%  "
$ws.Range("D34").Value = "'6.55"
$ws.Range("E34").Value = "  +4.24%  "
$ws.Range("D35").Value = "'156.17"
$ws.Range("E35").Value = "  -1.96%  "
$ws.Range("E36").Value = "  +0.98%  "
$ws.Range("D37").Value = "'2.829.82"
$ws.Range("E37").Value = "  +6.59%  "
$ws.Range("D38").Value = "'26.06"
$ws.Range("E38").Value = "  +3.49%  "
$ws.Range("D39").Value = "'0.0703"
$ws.Range("E39").Value = "  +2.37%  "
$ws.Range("D40").Value = "'1.67"
$ws.Range("E40").Value = "  -0.38%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "'0.0300"
$ws.Range("E41").Value = "  +4.95%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").Value = "'4.25"
$ws.Range("E42").Value = "  +0.14%  "
$ws.Range("D43").Value = "'39.81"
$ws.Range("E43").Value = "  +2.56%  "
$ws.Range("E44").Value = "  +0.75%  "
$ws.Range("E45").Value = "  +1.33%  "
$ws.Range("D46").Value = "'3.232.77"
$ws.Range("E46").Value = "  +1.13%  "
$ws.Range("D47").Value = "'0.985"
$ws.Range("E47").Value = "  -0.22%  "
$ws.Range("E48").Value = "  -0.66%  "
$ws.Range("D49").Value = "'20.68"
$ws.Range("E49").Value = "  +1.62%  "
$ws.Range("D50").Value = "'0.794"
$ws.Range("E50").Value = "  +4.84%  "
$ws.Range("D51").Value = "'1.00"
$ws.Range("E51").Value = "  -0.01%  "
